# "Fruta / hortaliza, semanal"
# Insert a new weekly record at row 35 (pushing the existing rows 35-95
# down to 36-96) and populate it with the new week's data. All the other
# rows keep their previous values, just shifted down by one row, which
# Excel's native row Insert() already takes care of.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 35:95 down to 36:96, leaving a blank row 35 behind.
$ws.Rows.Item(35).Insert()

# Fill the newly inserted row 35 with the new week's observation.
$ws.Range("A35").Value = 5
$ws.Range("B35").Value = "Macroferia Regional de Talca"
$ws.Range("C35").Value = "Maule"
$ws.Range("D35").Value = 44581
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 100112030
$ws.Range("G35").Value = "Poroto granado"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 25000
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = 25000
$ws.Range("N35").Value = "$/saco 25 kilos"
$ws.Range("O35").Value = "Región del Maule"
$ws.Range("P35").Value = 1000
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
